# Automatic update of files.
# Bump the "Förändrad" (Changed) date serial in column C from 46075 to 46076
# for every data row (rows 2 through 243).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 243

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46075) {
        $cell.Value = 46076
    }
}
